$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Numeric-looking Price strings are written with a leading apostrophe (quote-prefix)
# to keep them as text like the source data, then the original cell style is
# restored so no stray number formatting is left behind.

$ws.Range("D2").Value = '29.598.37'
$ws.Range("E2").Value = '  +3.36%  '

$ws.Range("D3").Value = '1.605.71'
$ws.Range("E3").Value = '  +2.70%  '

$ws.Range("E4").Value = '  -0.04%  '

$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'212.34"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +1.07%  '

$ws.Range("E6").Value = '  +2.56%  '

$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = '  -0.07%  '

$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'26.75"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = '  +7.39%  '

$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'43.57"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  -1.32%  '

$ws.Range("E10").Value = '  +2.24%  '

$ws.Range("E11").Value = '  +2.51%  '

$ws.Range("E12").Value = '  +1.27%  '

$ws.Range("D13").Value = '1.835.77'
$ws.Range("E13").Value = '  +2.66%  '

$ws.Range("D14").Value = '1.624.97'
$ws.Range("E14").Value = '  +3.88%  '

$ws.Range("D15").Value = '29.612.12'
$ws.Range("E15").Value = '  +3.24%  '

$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'0.536"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = '  +3.80%  '

$ws.Range("E17").Value = '  +2.30%  '

$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'63.44"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = '  +3.23%  '

$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'240.41"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  +5.67%  '

$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'7.60"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  +3.95%  '

$ws.Range("D21").Value = '0.0₃0692'
$ws.Range("E21").Value = '  +1.87%  '

$ws.Range("E22").Value = '  -0.04%  '

$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'3.99"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  +1.74%  '

$ws.Range("E24").Value = '  +2.17%  '

$ws.Range("E25").Value = '  +0.67%  '

$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'154.45"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  +1.88%  '

$ws.Range("E27").Value = '  +2.45%  '

$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'15.28"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  +3.48%  '

$ws.Range("E29").Value = '  +2.84%  '

$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("E31").Value = '  +3.33%  '

$ws.Range("E32").Value = '  +0.98%  '

$ws.Range("E33").Value = '  +1.60%  '

$ws.Range("E34").Value = '  +4.22%  '

$ws.Range("D35").Value = '1.407.41'
$ws.Range("E35").Value = '  +0.61%  '

$ws.Range("E36").Value = '  +0.49%  '

$ws.Range("E37").Value = '  +5.18%  '

$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'2.81"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  +5.67%  '

$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'2.31"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  +0.24%  '

$ws.Range("E40").Value = '  +2.58%  '

$ws.Range("E41").Value = '  +4.20%  '

$ws.Range("E42").Value = '  +2.58%  '

$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.0493"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  +7.20%  '

$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'54.02"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = '  +27.38%  '

$ws.Range("E45").Value = '  +3.83%  '

$ws.Range("E46").Value = '  -0.02%  '

$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'65.96"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = '  +3.34%  '

$ws.Range("E48").Value = '  +1.23%  '

$ws.Range("D49").Value = '1.746.14'
$ws.Range("E49").Value = '  +2.86%  '

$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'0.861"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  -0.24%  '

$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'86.67"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  +2.28%  '
